# Calculate.xlsx - "Test calc Arbitrage all 3 Factors"
#
# The sheet has two parallel "triangular arbitrage" blocks (rows 3-9 and
# rows 13-19). In the second block, the ETH/BTC and ETH/USDT quotes (and
# their prices) were entered in the wrong columns - this swaps C<->E for
# row 13 (headers) and row 14 (values) so the C16/E16 formulas evaluate
# the third arbitrage leg correctly, and mirrors the "Verkaufen"/"Ask"
# labels (A18:A19) into column C to label the new factor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the ETH/BTC / ETH/USDT headers in row 13.
$ws.Range("C13").Value = "ETH/USDT"
$ws.Range("E13").Value = "ETH/BTC"

# Swap the corresponding quote values in row 14.
$ws.Range("C14").Value = 1699.5
$ws.Range("E14").Value = 0.030954

# Add the third "Verkaufen" / "Ask" label column (mirrors A18:A19).
$ws.Range("C18").Value = "Verkaufen"
$ws.Range("C19").Value = "Ask"

# Move the selection to where the edit ended up.
$ws.Range("E18").Select()
